$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 2 per the FlashScore data refresh
$ws.Range("G2").Value = 1.44
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 1.95
$ws.Range("L2").Value = 6.5
$ws.Range("X2").Value = 7
$ws.Range("AB2").Value = 26
$ws.Range("AE2").Value = 19
$ws.Range("AL2").Value = 51
